$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated accuracy_all table (re-sorted by accuracy_condition, descending,
# with refreshed accuracy values) per commit "Commit 2 16th May".
$data = @(
    @("all alpha metrics", 0.878338278931751, 0.888724035608309),
    @("faith_pd + gini_index", 0.847181008902077, 0.859050445103858),
    @("menhinick + gini_index", 0.836795252225519, 0.851632047477745),
    @("chao1 + gini_index", 0.835311572700297, 0.86646884272997),
    @("fisher_alpha + gini_index", 0.829376854599407, 0.853115727002967),
    @("margalef + gini_index", 0.8234421364985161, 0.850148367952522),
    @("chao1 + pielou_evenness", 0.7195845697329381, 0.743323442136498),
    @("chao1 + simpson", 0.7195845697329381, 0.741839762611276),
    @("menhinick + strong", 0.716617210682493, 0.737388724035608),
    @("menhinick + pielou_evenness", 0.713649851632048, 0.732937685459941),
    @("menhinick + shannon_entropy", 0.70919881305638, 0.722551928783383),
    @("fisher_alpha + pielou_evenness", 0.70919881305638, 0.737388724035608),
    @("fisher_alpha + shannon_entropy", 0.707715133531157, 0.729970326409496),
    @("margalef + shannon_entropy", 0.704747774480712, 0.737388724035608),
    @("margalef + pielou_evenness", 0.6988130563798221, 0.734421364985163),
    @("margalef + simpson", 0.6988130563798221, 0.732937685459941),
    @("faith_pd + pielou_evenness", 0.694362017804154, 0.704747774480712),
    @("chao1 + shannon_entropy", 0.692878338278932, 0.72700296735905),
    @("fisher_alpha + simpson", 0.689910979228487, 0.706231454005935),
    @("menhinick + simpson", 0.682492581602374, 0.706231454005935),
    @("chao1 + strong", 0.6780415430267061, 0.724035608308605),
    @("margalef + strong", 0.6765578635014839, 0.712166172106825),
    @("faith_pd + shannon_entropy", 0.6765578635014839, 0.683976261127596),
    @("faith_pd + simpson", 0.675074183976261, 0.695845697329377),
    @("fisher_alpha + strong", 0.670623145400593, 0.713649851632048),
    @("faith_pd + strong", 0.642433234421365, 0.663204747774481)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
